$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue "D2" "243.43"
Set-TextValue "D3" "25.15"
Set-TextValue "D4" "5.162"
Set-TextValue "D5" "0.05721"
Set-TextValue "D6" "6.492"
Set-TextValue "D7" "3.111"
Set-TextValue "D8" "0.8097"
Set-TextValue "D9" "0.8459"
Set-TextValue "D10" "0.1338"
Set-TextValue "D11" "0.06949"
Set-TextValue "D12" "0.02832"
Set-TextValue "D13" "0.09369"
Set-TextValue "D14" "0.001512"
Set-TextValue "D15" "0.0005970"
$ws.Range("E15").Value = "14OneONE"
Set-TextValue "D16" "0.006223"
Set-TextValue "D17" "3.502"
Set-TextValue "D19" "0.3199"
Set-TextValue "D20" "0.03152"
Set-TextValue "D21" "0.1301"
Set-TextValue "D22" "3.751"
Set-TextValue "D23" "0.04656"
Set-TextValue "D24" "0.1329"
Set-TextValue "D25" "0.001234"
Set-TextValue "D26" "0.004259"
Set-TextValue "D27" "0.00009699"
$ws.Range("E27").Value = "26NitroExNTX"
Set-TextValue "D40" "0.03617"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
Set-TextValue "D43" "0.003000"
Set-TextValue "D44" "0.007391"
Set-TextValue "D45" "0.00005286"
Set-TextValue "D47" "0.2100"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
